# Sort alternate source field columns on sheet "Forms1":
# Swap the content of columns D and E (header row 1 and data row 3),
# so that "Alternate Source Field 1"/"sex" and
# "Alternate Source Field 2"/"sex0" trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forms1")

# Row 1 header swap: D1 <-> E1
$d1 = $ws.Range("D1").Value2
$e1 = $ws.Range("E1").Value2
$ws.Range("D1").Value = $e1
$ws.Range("E1").Value = $d1

# Row 3 data swap: D3 <-> E3
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2
$ws.Range("D3").Value = $e3
$ws.Range("E3").Value = $d3

# Update the selection to reflect the post-edit active cell
$ws.Range("D1").Select()
